# Edit script: applies the diff to ANEXO A Practicas.docx
# Strategy: pull the Flat-OPC WordOpenXML string, perform the same
# literal XML edits the diff describes (text resplits with proofErr
# marks, bookmark move, comment id renumbering), and write it back.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# ---- word/document.xml edits ----
$old = '<w:ind w:right="6"/><w:jc w:val="center"/><w:rPr><w:bCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r w:rsidRPr="007E184D"><w:rPr><w:bCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>COMPLETO</w:t></w:r></w:p>'
$new = '<w:ind w:right="6"/><w:rPr><w:bCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r w:rsidRPr="007E184D"><w:rPr><w:bCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>COMPLETO</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #0" }
$xml = $xml.Replace($old, $new)

$old = '<w:r w:rsidRPr="00C85302"><w:rPr><w:bCs/></w:rPr><w:t>Analizar forma de hacer camabios en la variante operativa dado una situación determina de Urgencia</w:t></w:r>'
$new = '<w:r w:rsidRPr="00C85302"><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">Analizar forma de hacer </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>camabios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> en la variante operativa dado una situación determina de Urgencia</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #1" }
$xml = $xml.Replace($old, $new)

$old = '<w:r w:rsidRPr="00C85302"><w:rPr><w:bCs/></w:rPr><w:t>Hacer un estudio de cómo se hacen las cofiguraciones a partir de una emergencia y el uso de la simulación para dichas configuraciones</w:t></w:r>'
$new = '<w:r w:rsidRPr="00C85302"><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">Hacer un estudio de cómo se hacen las </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>cofiguraciones</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> a partir de una emergencia y el uso de la simulación para dichas configuraciones</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #2" }
$xml = $xml.Replace($old, $new)

$old = '<w:r w:rsidRPr="00C85302"><w:rPr><w:bCs/></w:rPr><w:t>Analizar la propuesta de la tesis de Magaret Sánchez</w:t></w:r>'
$new = '<w:r w:rsidRPr="00C85302"><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">Analizar la propuesta de la tesis de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>Magaret</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> Sánchez</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #3" }
$xml = $xml.Replace($old, $new)

$old = '<w:r w:rsidRPr="00C85302"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CU"/></w:rPr><w:t>Ubicación Inicial de estaciones de servicio (dígase a apartir del MCLPTools, o por configuración de usuario)</w:t></w:r>'
$new = '<w:r w:rsidRPr="00C85302"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CU"/></w:rPr><w:t xml:space="preserve">Ubicación Inicial de estaciones de servicio (dígase a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CU"/></w:rPr><w:t>apartir</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CU"/></w:rPr><w:t xml:space="preserve"> del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CU"/></w:rPr><w:t>MCLPTools</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CU"/></w:rPr><w:t>, o por configuración de usuario)</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #4" }
$xml = $xml.Replace($old, $new)

$old = '<w:r><w:rPr><w:bCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>14</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$new = '<w:r><w:rPr><w:bCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>14</w:t></w:r></w:p>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #5" }
$xml = $xml.Replace($old, $new)

$old = '<w:r w:rsidRPr="001C4DC2"><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>-Jefe de proyecto</w:t></w:r>'
$new = '<w:r w:rsidRPr="001C4DC2"><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Jefe</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> de proyecto</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #6" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeStart w:id="7"/>'
$new = '<w:commentRangeStart w:id="8"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #7" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeEnd w:id="7"/>'
$new = '<w:commentRangeEnd w:id="8"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #8" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentReference w:id="7"/>'
$new = '<w:commentReference w:id="8"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #9" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeStart w:id="6"/>'
$new = '<w:commentRangeStart w:id="7"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #10" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeEnd w:id="6"/>'
$new = '<w:commentRangeEnd w:id="7"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #11" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentReference w:id="6"/>'
$new = '<w:commentReference w:id="7"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #12" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeStart w:id="5"/>'
$new = '<w:commentRangeStart w:id="6"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #13" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeEnd w:id="5"/>'
$new = '<w:commentRangeEnd w:id="6"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #14" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentReference w:id="5"/>'
$new = '<w:commentReference w:id="6"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #15" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeStart w:id="4"/>'
$new = '<w:commentRangeStart w:id="5"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #16" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeEnd w:id="4"/>'
$new = '<w:commentRangeEnd w:id="5"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #17" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentReference w:id="4"/>'
$new = '<w:commentReference w:id="5"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #18" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeStart w:id="3"/>'
$new = '<w:commentRangeStart w:id="4"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #19" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeEnd w:id="3"/>'
$new = '<w:commentRangeEnd w:id="4"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #20" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentReference w:id="3"/>'
$new = '<w:commentReference w:id="4"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #21" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeStart w:id="2"/>'
$new = '<w:commentRangeStart w:id="3"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #22" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeEnd w:id="2"/>'
$new = '<w:commentRangeEnd w:id="3"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #23" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentReference w:id="2"/>'
$new = '<w:commentReference w:id="3"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #24" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeStart w:id="1"/>'
$new = '<w:commentRangeStart w:id="2"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #25" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeEnd w:id="1"/>'
$new = '<w:commentRangeEnd w:id="2"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #26" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentReference w:id="1"/>'
$new = '<w:commentReference w:id="2"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #27" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeStart w:id="0"/>'
$new = '<w:commentRangeStart w:id="1"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #28" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentRangeEnd w:id="0"/>'
$new = '<w:commentRangeEnd w:id="1"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #29" }
$xml = $xml.Replace($old, $new)

$old = '<w:commentReference w:id="0"/>'
$new = '<w:commentReference w:id="1"/>'
if ($xml.IndexOf($old) -lt 0) { throw "document.xml pattern not found for edit #30" }
$xml = $xml.Replace($old, $new)

# ---- word/comments.xml edits ----
$old = '<w:r><w:t>Básicamente hay que diseñar una arquitectura que interactúe con la capa de lógica del framework para que los agentes interactúen entre si. Y Esta arquitectura va a interactuar con la capa controlador de la arquitectura Modelo-vista-controlador con la que viene MASON</w:t></w:r>'
$new = '<w:r><w:t xml:space="preserve">Básicamente hay que diseñar una arquitectura que interactúe con la capa de lógica del framework para que los agentes interactúen entre </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>si</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. Y Esta arquitectura va a interactuar con la capa controlador de la arquitectura Modelo-vista-controlador con la que viene MASON</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #0" }
$xml = $xml.Replace($old, $new)

$old = '<w:r><w:t>Ya esto esta echo, al menos una version inicial, todavía la estoy retocando xq esta enredado, los otros días publique una duda sobre como hacerlo y nadie ha hecho nada parecido, xq nadie trabaja con tantos agentes como nosotros, asi que estoy leyendo muchos libros de ingeniería de softweare para hacer algo que sirva y no una caquita</w:t></w:r>'
$new = '<w:r><w:t xml:space="preserve">Ya esto </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>esta</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> echo, al menos una </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>version</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> inicial, todavía la estoy retocando </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>xq</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> esta enredado, los otros días publique una duda sobre </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>como</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> hacerlo y nadie ha hecho nada parecido, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>xq</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> nadie trabaja con tantos agentes como nosotros, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>asi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> que estoy leyendo muchos libros de ingeniería de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>softweare</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para hacer algo que sirva y no una caquita</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #1" }
$xml = $xml.Replace($old, $new)

$old = '<w:r><w:t>Básicamente lo mismo, pero para la parte visual, xq MASON tiene MUY separada estas capas y hay que hacerlo básicamente 2veces</w:t></w:r>'
$new = '<w:r><w:t xml:space="preserve">Básicamente lo mismo, pero para la parte visual, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>xq</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> MASON tiene MUY separada estas capas y hay que hacerlo básicamente 2veces</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #2" }
$xml = $xml.Replace($old, $new)

$old = '<w:r><w:t>Aquí es donde va todos los cambios visuales que le queremos hacer, poner iconitos chulos y esas cosas, esto va de la mano con la parte de la arquitectura visual de los modulos, ya que cada modulo es quien define la visualización que va a tener</w:t></w:r>'
$new = '<w:r><w:t xml:space="preserve">Aquí es donde va todos los cambios visuales que le queremos hacer, poner iconitos chulos y esas cosas, esto va de la mano con la parte de la arquitectura visual de los </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>modulos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, ya que cada </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>modulo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> es quien define la visualización que va a tener</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #3" }
$xml = $xml.Replace($old, $new)

$old = '<w:r><w:t>Este es el modelo tal y como está ahora, no le voy a agregar mas nada para que me de tiempo a hacer todo lo demás y dejar el sistema bien echo para cuando empecemos la tesis poder agregarle mas cosas, todos los agentes que queramos y los comportamientos y demás. Claro, lo que, si hay que retocarlo y hacerlo bien, porque hora mismo esta al berro.</w:t></w:r>'
$new = '<w:r><w:t xml:space="preserve">Este es el modelo tal y como está ahora, no le voy a agregar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> nada para que me </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>de</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tiempo a hacer todo lo demás y dejar el sistema bien echo para cuando empecemos la tesis poder agregarle </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cosas, todos los agentes que queramos y los comportamientos y demás. Claro, lo que, si hay que retocarlo y hacerlo bien, porque hora mismo esta al berro.</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #4" }
$xml = $xml.Replace($old, $new)

$old = '<w:r><w:t>Aquí en realidad quiero agregarle al menos bomberos muy básicos para probar lo de cargar las ubicaciones por ficheros y eso, y que 2 tipos de agentes atiendan la misma emergencia, pero no lo pongo especifico asi no me vaya a quedar corto de tiempo y después me falte algo por hacer</w:t></w:r>'
$new = '<w:r><w:t xml:space="preserve">Aquí en realidad quiero agregarle al menos bomberos muy básicos para probar lo de cargar las ubicaciones por ficheros y eso, y que 2 tipos de agentes atiendan la misma emergencia, pero no lo pongo especifico </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>asi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> no me vaya a quedar corto de tiempo y después me falte algo por hacer</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #5" }
$xml = $xml.Replace($old, $new)

$old = '<w:r><w:t>Necesito que me digas aquí, que tiene que tener el informe, xq el que yo entregue en 3ro, era el mismo de margaret de 4to, que hago cojo ese y le agrego lo que te falta, o se lo agrego ya directo a la tesis</w:t></w:r>'
$new = '<w:r><w:t xml:space="preserve">Necesito que me digas aquí, que tiene que tener el informe, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>xq</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> el que yo entregue en 3ro, era el mismo de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>margaret</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de 4to, que hago cojo ese y le agrego lo que te falta, o se lo agrego ya directo a la tesis</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #6" }
$xml = $xml.Replace($old, $new)

$old = '<w:r><w:t>Aquí estaba el dia 13, no tiene sentido xq tiene que estar en el rango [15,21], le puse 18 que es aproximadamente el del medio</w:t></w:r>'
$new = '<w:r><w:t xml:space="preserve">Aquí estaba el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 13, no tiene sentido </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>xq</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tiene que estar en el rango [15,21], le puse 18 que es aproximadamente el del medio</w:t></w:r>'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #7" }
$xml = $xml.Replace($old, $new)

$old = '<w:comment w:id="7" w:author="Yo" w:date="2020-11-10T20:23:00Z" w:initials="Y">'
$new = '<w:comment w:id="8" w:author="Yo" w:date="2020-11-10T20:23:00Z" w:initials="Y">'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #8" }
$xml = $xml.Replace($old, $new)

$old = '<w:comment w:id="6" w:author="Yo" w:date="2020-11-10T21:01:00Z" w:initials="Y">'
$new = '<w:comment w:id="7" w:author="Yo" w:date="2020-11-10T21:01:00Z" w:initials="Y">'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #9" }
$xml = $xml.Replace($old, $new)

$old = '<w:comment w:id="5" w:author="Yo" w:date="2020-11-10T20:55:00Z" w:initials="Y">'
$new = '<w:comment w:id="6" w:author="Yo" w:date="2020-11-10T20:55:00Z" w:initials="Y">'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #10" }
$xml = $xml.Replace($old, $new)

$old = '<w:comment w:id="4" w:author="Yo" w:date="2020-11-10T20:49:00Z" w:initials="Y">'
$new = '<w:comment w:id="5" w:author="Yo" w:date="2020-11-10T20:49:00Z" w:initials="Y">'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #11" }
$xml = $xml.Replace($old, $new)

$old = '<w:comment w:id="3" w:author="Yo" w:date="2020-11-10T20:42:00Z" w:initials="Y">'
$new = '<w:comment w:id="4" w:author="Yo" w:date="2020-11-10T20:42:00Z" w:initials="Y">'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #12" }
$xml = $xml.Replace($old, $new)

$old = '<w:comment w:id="2" w:author="Yo" w:date="2020-11-10T20:37:00Z" w:initials="Y">'
$new = '<w:comment w:id="3" w:author="Yo" w:date="2020-11-10T20:37:00Z" w:initials="Y">'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #13" }
$xml = $xml.Replace($old, $new)

$old = '<w:comment w:id="1" w:author="Yo" w:date="2020-11-10T20:34:00Z" w:initials="Y">'
$new = '<w:comment w:id="2" w:author="Yo" w:date="2020-11-10T20:34:00Z" w:initials="Y">'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #14" }
$xml = $xml.Replace($old, $new)

$old = '<w:comment w:id="0" w:author="Yo" w:date="2020-11-10T20:39:00Z" w:initials="Y">'
$new = '<w:comment w:id="1" w:author="Yo" w:date="2020-11-10T20:39:00Z" w:initials="Y">'
if ($xml.IndexOf($old) -lt 0) { throw "comments.xml pattern not found for edit #15" }
$xml = $xml.Replace($old, $new)

$d.WordOpenXML = $xml
Write-Output "Edit applied successfully"